# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.692.05"
$ws.Range("E2").Value = "  -1.55%  "

# Row 3
$ws.Range("D3").Value = "2.443.68"
$ws.Range("E3").Value = "  -1.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.39%  "

# Row 5
$ws.Range("D5").Value = "'569.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "

# Row 6
$ws.Range("D6").Value = "'145.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.97%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("E8").Value = "  -1.55%  "

# Row 9
$ws.Range("D9").Value = "2.437.28"
$ws.Range("E9").Value = "  -2.73%  "

# Row 10
$ws.Range("E10").Value = "  -4.32%  "

# Row 11
$ws.Range("D11").Value = "'0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.60%  "

# Row 12
$ws.Range("D12").Value = "'5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.70%  "

# Row 13
$ws.Range("E13").Value = "  -2.63%  "

# Row 14
$ws.Range("D14").Value = "'27.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "

# Row 15
$ws.Range("D15").Value = "'0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.21%  "

# Row 16
$ws.Range("E16").Value = "  -1.57%  "

# Row 17
$ws.Range("D17").Value = "62.559.26"
$ws.Range("E17").Value = "  -1.75%  "

# Row 18
$ws.Range("D18").Value = "2.437.43"
$ws.Range("E18").Value = "  -2.51%  "

# Row 19
$ws.Range("D19").Value = "'11.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.88%  "

# Row 20
$ws.Range("D20").Value = "'7.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "

# Row 21
$ws.Range("D21").Value = "'327.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "

# Row 22
$ws.Range("E22").Value = "  -2.11%  "

# Row 23
$ws.Range("D23").Value = "'2.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.80%  "

# Row 24
$ws.Range("E24").Value = "  -0.54%  "

# Row 25
$ws.Range("D25").Value = "'65.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.51%  "

# Row 26
$ws.Range("D26").Value = "'620.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.26%  "

# Row 27
$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0998"
$ws.Range("E28").Value = "  -5.52%  "

# Row 29
$ws.Range("D29").Value = "2.556.25"
$ws.Range("E29").Value = "  -1.73%  "

# Row 30
$ws.Range("E30").Value = "  -1.85%  "

# Row 31
$ws.Range("E31").Value = "  +0.39%  "

# Row 32
$ws.Range("E32").Value = "  -5.25%  "

# Row 33
$ws.Range("D33").Value = "'1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.25%  "

# Row 34
$ws.Range("E34").Value = "  -4.55%  "

# Row 35
$ws.Range("D35").Value = "'5.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "

# Row 36
$ws.Range("E36").Value = "  -3.92%  "

# Row 37
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38
$ws.Range("D38").Value = "'0.377"
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'18.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.27%  "

# Row 40
$ws.Range("E40").Value = "  -4.46%  "

# Row 41
$ws.Range("D41").Value = "'146.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.94%  "

# Row 42
$ws.Range("E42").Value = "  -4.80%  "

# Row 43
$ws.Range("D43").Value = "'2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.76%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").Value = "'20.68"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.0528"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.21%  "

# Row 49
$ws.Range("E49").Value = "  -3.06%  "

# Row 50
$ws.Range("E50").Value = "  -3.17%  "

# Row 51
$ws.Range("D51").Value = "'0.0917"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.53%  "

# Rows 45/46: Filecoin and Aave swapped ranking order
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'3.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'146.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.69%  "

